# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at the top of the Alcachofa data
# block (row 118), pushing the existing rows 118-138 down to 119-139.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 118 (default Insert() on a full row shifts
# the existing rows down, i.e. xlShiftDown), carrying over formatting
# (date style) from the row above, same as a manual Excel "Insert Row".
$ws.Rows("118:118").Insert()

# Populate the newly inserted row with the new data point.
$ws.Cells.Item(118, 1).Value = 7
$ws.Cells.Item(118, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(118, 3).Value = 'Ñuble'
$ws.Cells.Item(118, 4).Value = 45218
$ws.Cells.Item(118, 5).Value = 16
$ws.Cells.Item(118, 6).Value = 100112013
$ws.Cells.Item(118, 7).Value = 'Alcachofa'
$ws.Cells.Item(118, 8).Value = 'Española'
$ws.Cells.Item(118, 9).Value = 'Primera'
$ws.Cells.Item(118, 10).Value = 80
$ws.Cells.Item(118, 11).Value = 13000
$ws.Cells.Item(118, 12).Value = 13000
$ws.Cells.Item(118, 13).Value = 13000
$ws.Cells.Item(118, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(118, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(118, 16).Value = 433
$ws.Cells.Item(118, 17).Value = 30
$ws.Cells.Item(118, 18).Value = 'Hortaliza'
